# Adding get registration pin from mail subject functionality
# Adds two new LinkedIn user rows (Ahmed Elemam / Islam Azez) to the data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 5: Ahmed Elemam ----
$ws.Range("A5").Value = "Ahmed"
$ws.Range("B5").Value = "Elemam"
$ws.Range("C5").Value = "aemamarcom83@gmail.com"
$ws.Range("D5").Value = "P@ssw0rd85"
$ws.Range("E5").Value = "Italy"
$ws.Range("F5").Value = 155147
$ws.Range("G5").Value = "CEO"
$ws.Range("H5").Value = "Alqemam"
$ws.Range("I5").Value = "Hospital & Health Care"

# ---- Row 6: Islam Azez ----
$ws.Range("A6").Value = "Islam"
$ws.Range("B6").Value = "Azez"
$ws.Range("C6").Value = "islamazez83@gmail.com"
$ws.Range("D6").Value = "P@ssw0rd85"
$ws.Range("E6").Value = "Egypt"
$ws.Range("F6").Value = 35111
$ws.Range("G6").Value = "Senior Software QC Engineer"
$ws.Range("H6").Value = "ARCOM"
$ws.Range("I6").Value = "Information Technology and Services"

# ---- mailto: hyperlinks, created in the same order the workbook records them ----
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:aemamarcom83@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:P@ssw0rd85")
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:P@ssw0rd85")
$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:islamazez83@gmail.com")

# Re-apply the sheet's "Hyperlink" style to the e-mail/password cells (same look as rows 2-4)
$ws.Range("C2").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("D5").PasteSpecial(-4122)
$ws.Range("D6").PasteSpecial(-4122)
$ws.Range("C6").PasteSpecial(-4122)

# I6 (Work Industry) reuses the small Consolas "pin" look already used on I4
$ws.Range("I4").Copy()
$ws.Range("I6").PasteSpecial(-4122)

# ---- Column G needs to grow to fit "Senior Software QC Engineer" ----
$ws.Columns("G").ColumnWidth = 26.166666666666668

# ---- Selection ends on the last cell touched ----
$ws.Range("F6").Select()
